$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the contents of columns A, B, E, F, G, H, Q, R
# between row 2 and row 4 (all other columns are identical between
# the two rows, so only these need to change).
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range("$col" + "2")
    $cell4 = $ws.Range("$col" + "4")
    $tmp = $cell2.Value2
    $cell2.Value2 = $cell4.Value2
    $cell4.Value2 = $tmp
}
